$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("B1").Value = "e_modulus"
$ws.Range("C1").Value = "tensile_yield_strength"
$ws.Range("D1").Value = "tensile_strain_at_break"

# Update row labels (column A)
$ws.Range("A2").Value = "e_modulus"
$ws.Range("A3").Value = "tensile_yield_strength"
$ws.Range("A4").Value = "tensile_strain_at_break"

# Update selection
$ws.Range("D6").Select()
